# fix(publipostage): Try to solve Excel emoji problem
#
# The "statut" column (A) used emoji characters as status codes. Excel had
# trouble handling these emoji, so replace them with plain-text equivalents:
#   📕 -> "-3"
#   📘 -> "⚠️"
#   📙 -> "+3"
#   📗 -> "✅"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

# Column A holds the "statut" values (row 1 is the header "statut").
$col = 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value()
    if ($val -isnot [string]) {
        continue
    }

    if ($val -eq "📕") {
        # "-3" looks numeric, force the cell to stay text so it round-trips
        # as a string instead of being coerced into a number.
        $cell.NumberFormat = "@"
        $cell.Value = "-3"
    } elseif ($val -eq "📘") {
        $cell.Value = "⚠️"
    } elseif ($val -eq "📙") {
        $cell.NumberFormat = "@"
        $cell.Value = "+3"
    } elseif ($val -eq "📗") {
        $cell.Value = "✅"
    }
}
